$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3406.25
$ws.Range("J74").Value = 3406.25
$ws.Range("L74").Value = 3406.25
$ws.Range("N74").Value = -5278.25

$ws.Range("H77").Value = 3406.25
$ws.Range("J77").Value = 3406.25
$ws.Range("L77").Value = 17031.25
$ws.Range("N77").Value = -26391.25

$ws.Range("H86").Value = 1894.8948
$ws.Range("I86").Value = 1813.5333
$ws.Range("J86").Value = 2200
$ws.Range("K86").Value = 1813.5333
$ws.Range("L86").Value = 2200
$ws.Range("M86").Value = -690.5333000000001
$ws.Range("N86").Value = -4446

$ws.Range("H88").Value = 7556
$ws.Range("I88").Value = 960.3333
$ws.Range("J88").Value = 14151.667
$ws.Range("K88").Value = 960.3333
$ws.Range("L88").Value = 14151.667
$ws.Range("M88").Value = -554.3333
$ws.Range("N88").Value = -14963.667

$ws.Range("H89").Value = 1894.8948
$ws.Range("I89").Value = 1813.5333
$ws.Range("J89").Value = 2200
$ws.Range("K89").Value = 9067.666499999999
$ws.Range("L89").Value = 11000
$ws.Range("M89").Value = -3451.666499999999
$ws.Range("N89").Value = -22232

$ws.Range("H91").Value = 7556
$ws.Range("I91").Value = 960.3333
$ws.Range("J91").Value = 14151.667
$ws.Range("K91").Value = 960.3333
$ws.Range("L91").Value = 14151.667
$ws.Range("M91").Value = 443.6667
$ws.Range("N91").Value = -16959.667

$ws.Range("H132").Value = 2268.9714
$ws.Range("I132").Value = 2250.4666
$ws.Range("J132").Value = 2380
$ws.Range("K132").Value = 6751.399800000001
$ws.Range("L132").Value = 7140
$ws.Range("M132").Value = -4221.399800000001
$ws.Range("N132").Value = -12200

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1625.4286
$ws.Range("I61").Value = 1228.5
$ws.Range("J61").Value = 4007
$ws.Range("K61").Value = 1228.5
$ws.Range("L61").Value = 4007
$ws.Range("M61").Value = -1016.5
$ws.Range("N61").Value = -4431

$ws.Range("H136").Value = 1625.4286
$ws.Range("I136").Value = 1228.5
$ws.Range("J136").Value = 4007
$ws.Range("K136").Value = 3685.5
$ws.Range("L136").Value = 12021
$ws.Range("M136").Value = -1135.5
$ws.Range("N136").Value = -17121

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 12942.833
$ws.Range("I75").Value = 12942.833
$ws.Range("K75").Value = 12942.833
$ws.Range("M75").Value = -12006.833

$ws.Range("H78").Value = 12942.833
$ws.Range("I78").Value = 12942.833
$ws.Range("K78").Value = 38828.499
$ws.Range("M78").Value = -34148.499

$ws.Range("H82").Value = 13160.5
$ws.Range("I82").Value = 9620
$ws.Range("J82").Value = 20241.5
$ws.Range("K82").Value = 9620
$ws.Range("L82").Value = 20241.5
$ws.Range("M82").Value = -9237
$ws.Range("N82").Value = -21007.5

$ws.Range("H85").Value = 13160.5
$ws.Range("I85").Value = 9620
$ws.Range("J85").Value = 20241.5
$ws.Range("K85").Value = 9620
$ws.Range("L85").Value = 20241.5
$ws.Range("M85").Value = -8294
$ws.Range("N85").Value = -22893.5

$ws.Range("H94").Value = 861.73334
$ws.Range("I94").Value = 840.46155
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 840.46155
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -389.46155
$ws.Range("N94").Value = -1902

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2160.3635
$ws.Range("I58").Value = 1533.7142
$ws.Range("J58").Value = 3257
$ws.Range("K58").Value = 1533.7142
$ws.Range("L58").Value = 3257
$ws.Range("M58").Value = -1330.7142
$ws.Range("N58").Value = -3663

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H136").Value = 2160.3635
$ws.Range("I136").Value = 1533.7142
$ws.Range("J136").Value = 3257
$ws.Range("K136").Value = 4601.142599999999
$ws.Range("L136").Value = 9771
$ws.Range("M136").Value = -2051.142599999999
$ws.Range("N136").Value = -14871

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 387.5
$ws.Range("I97").Value = 380.2
$ws.Range("J97").Value = 394.8
$ws.Range("K97").Value = 1140.6
$ws.Range("L97").Value = 1184.4
$ws.Range("M97").Value = -644.5999999999999
$ws.Range("N97").Value = -2176.4

$ws.Range("H132").Value = 1011692.7
$ws.Range("I132").Value = 1228
$ws.Range("K132").Value = 11052
$ws.Range("M132").Value = -8522

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8200
$ws.Range("I70").Value = 8200
$ws.Range("K70").Value = 8200
$ws.Range("M70").Value = -7930

$ws.Range("H73").Value = 8200
$ws.Range("I73").Value = 8200
$ws.Range("K73").Value = 8200
$ws.Range("M73").Value = -7264

$ws.Range("H80").Value = 2990.0688
$ws.Range("I80").Value = 2825
$ws.Range("K80").Value = 2825
$ws.Range("M80").Value = -1827

$ws.Range("H83").Value = 2990.0688
$ws.Range("I83").Value = 2825
$ws.Range("K83").Value = 14125
$ws.Range("M83").Value = -9133

$ws.Range("H122").Value = 6668589
$ws.Range("I122").Value = 12501650
$ws.Range("J122").Value = 2234.1428
$ws.Range("K122").Value = 37504950
$ws.Range("L122").Value = 6702.428400000001
$ws.Range("M122").Value = -37502500
$ws.Range("N122").Value = -11602.4284

$ws.Range("H132").Value = 2712.2083
$ws.Range("I132").Value = 2391.4211
$ws.Range("J132").Value = 3931.2
$ws.Range("K132").Value = 7174.263300000001
$ws.Range("L132").Value = 11793.6
$ws.Range("M132").Value = -4644.263300000001
$ws.Range("N132").Value = -16853.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 440
$ws.Range("I22").Value = 366.66666
$ws.Range("J22").Value = 550
$ws.Range("K22").Value = 366.66666
$ws.Range("L22").Value = 550
$ws.Range("M22").Value = -71.66665999999998
$ws.Range("N22").Value = -1140

$ws.Range("H27").Value = 440
$ws.Range("I27").Value = 366.66666
$ws.Range("J27").Value = 550
$ws.Range("K27").Value = 366.66666
$ws.Range("L27").Value = 550
$ws.Range("M27").Value = -259.66666
$ws.Range("N27").Value = -764

$ws.Range("H32").Value = 578.25
$ws.Range("I32").Value = 578.25
$ws.Range("K32").Value = 578.25
$ws.Range("M32").Value = -261.25

$ws.Range("H40").Value = 5992.5
$ws.Range("I40").Value = 8280.799999999999
$ws.Range("J40").Value = 4952.364
$ws.Range("K40").Value = 8280.799999999999
$ws.Range("L40").Value = 4952.364
$ws.Range("M40").Value = -8144.799999999999
$ws.Range("N40").Value = -5224.364

$ws.Range("H82").Value = 760
$ws.Range("I82").Value = 650
$ws.Range("J82").Value = 833.3333
$ws.Range("K82").Value = 650
$ws.Range("L82").Value = 833.3333
$ws.Range("M82").Value = -289
$ws.Range("N82").Value = -1555.3333

$ws.Range("H85").Value = 760
$ws.Range("I85").Value = 650
$ws.Range("J85").Value = 833.3333
$ws.Range("K85").Value = 650
$ws.Range("L85").Value = 833.3333
$ws.Range("M85").Value = 598
$ws.Range("N85").Value = -3329.3333

$ws.Range("H93").Value = 8706.916999999999
$ws.Range("I93").Value = 9453
$ws.Range("J93").Value = 500
$ws.Range("K93").Value = 9453
$ws.Range("L93").Value = 500
$ws.Range("M93").Value = -8205
$ws.Range("N93").Value = -2996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H122").Value = 1683.6578
$ws.Range("I122").Value = 1217.84
$ws.Range("K122").Value = 3653.52
$ws.Range("M122").Value = -1203.52
